$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Municipio" -> "CVE_MUN"
$ws.Range("A1").Value = "CVE_MUN"

# Format column A (rows 2-85) as Text so the numeric-looking municipality
# codes (e.g. "13001") are stored as strings, not numbers.
$ws.Range("A2:A85").NumberFormat = "@"

# Replace each municipality name with its corresponding CVE_MUN code,
# in the same row order as the original municipality list.
$ws.Range("A2").Value = "13001"
$ws.Range("A3").Value = "13002"
$ws.Range("A4").Value = "13003"
$ws.Range("A5").Value = "13004"
$ws.Range("A6").Value = "13005"
$ws.Range("A7").Value = "13006"
$ws.Range("A8").Value = "13007"
$ws.Range("A9").Value = "13008"
$ws.Range("A10").Value = "13010"
$ws.Range("A11").Value = "13011"
$ws.Range("A12").Value = "13013"
$ws.Range("A13").Value = "13012"
$ws.Range("A14").Value = "13014"
$ws.Range("A15").Value = "13015"
$ws.Range("A16").Value = "13017"
$ws.Range("A17").Value = "13018"
$ws.Range("A18").Value = "13019"
$ws.Range("A19").Value = "13016"
$ws.Range("A20").Value = "13009"
$ws.Range("A21").Value = "13020"
$ws.Range("A22").Value = "13021"
$ws.Range("A23").Value = "13022"
$ws.Range("A24").Value = "13023"
$ws.Range("A25").Value = "13024"
$ws.Range("A26").Value = "13025"
$ws.Range("A27").Value = "13026"
$ws.Range("A28").Value = "13027"
$ws.Range("A29").Value = "13028"
$ws.Range("A30").Value = "13029"
$ws.Range("A31").Value = "13030"
$ws.Range("A32").Value = "13031"
$ws.Range("A33").Value = "13032"
$ws.Range("A34").Value = "13033"
$ws.Range("A35").Value = "13040"
$ws.Range("A36").Value = "13034"
$ws.Range("A37").Value = "13035"
$ws.Range("A38").Value = "13037"
$ws.Range("A39").Value = "13051"
$ws.Range("A40").Value = "13038"
$ws.Range("A41").Value = "13039"
$ws.Range("A42").Value = "13041"
$ws.Range("A43").Value = "13042"
$ws.Range("A44").Value = "13043"
$ws.Range("A45").Value = "13044"
$ws.Range("A46").Value = "13045"
$ws.Range("A47").Value = "13048"
$ws.Range("A48").Value = "13047"
$ws.Range("A49").Value = "13049"
$ws.Range("A50").Value = "13050"
$ws.Range("A51").Value = "13036"
$ws.Range("A52").Value = "13052"
$ws.Range("A53").Value = "13053"
$ws.Range("A54").Value = "13046"
$ws.Range("A55").Value = "13054"
$ws.Range("A56").Value = "13055"
$ws.Range("A57").Value = "13056"
$ws.Range("A58").Value = "13057"
$ws.Range("A59").Value = "13058"
$ws.Range("A60").Value = "13059"
$ws.Range("A61").Value = "13060"
$ws.Range("A62").Value = "13061"
$ws.Range("A63").Value = "13062"
$ws.Range("A64").Value = "13063"
$ws.Range("A65").Value = "13064"
$ws.Range("A66").Value = "13065"
$ws.Range("A67").Value = "13067"
$ws.Range("A68").Value = "13068"
$ws.Range("A69").Value = "13069"
$ws.Range("A70").Value = "13070"
$ws.Range("A71").Value = "13071"
$ws.Range("A72").Value = "13072"
$ws.Range("A73").Value = "13073"
$ws.Range("A74").Value = "13074"
$ws.Range("A75").Value = "13075"
$ws.Range("A76").Value = "13076"
$ws.Range("A77").Value = "13077"
$ws.Range("A78").Value = "13066"
$ws.Range("A79").Value = "13078"
$ws.Range("A80").Value = "13079"
$ws.Range("A81").Value = "13080"
$ws.Range("A82").Value = "13081"
$ws.Range("A83").Value = "13082"
$ws.Range("A84").Value = "13083"
$ws.Range("A85").Value = "13084"
